# TimeLog_ConnorPeper.xlsx update:
#  - Append a sentence about the joinSession bug fix to the week's activity
#    log note (cell F12).
#  - Update the hours worked for that same week from 3.5 to 5 (cell E12).
#  - Move the view/selection (window scrolled to C4, active cell E13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Activities note (F12): append the bug-fix sentence to the existing text ---
$activityCell = $ws.Range("F12")
$existingNote = $activityCell.Text
$activityCell.Value = $existingNote + " Fixed the bug that had the session ID always be 2 if the session did not exist."

# --- Hours worked (E12): 3.5 -> 5 ---
$ws.Range("E12").Value = 5

# --- View state: scroll window top-left to C4 and select E13 ---
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("C4")
$ws.Range("E13").Select()
